# Applies:
#  1. The cached "datetimeFigureOut" auto-date field text, 12/13/19 -> 2/6/20,
#     on the slide master, every slide layout, and the notes master
#     (mirrors PowerPoint re-stamping the "Update automatically" date/time
#     placeholder the next time the deck is touched).
#  2. The FIFO-trend schematic's bivalve-density callout on slide 1:
#     "100 bivalves per foot" -> "400 cm of bivalves per foot".

$p = $ppt.ActivePresentation

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*" -and $sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq "12/13/19") {
                $sh.TextFrame.TextRange.Text = "2/6/20"
            }
        }
    }
}

# Slide master.
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every slide layout off the master.
for ($L = 1; $L -le $p.SlideMaster.CustomLayouts.Count; $L++) {
    $layout = $p.SlideMaster.CustomLayouts.Item($L)
    Update-DatePlaceholder $layout.Shapes
}

# Notes master.
Update-DatePlaceholder $p.NotesMaster.Shapes

# Slide 1: bivalve-density callout ("TextBox 42").
$s = $p.Slides.Item(1)
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -eq "TextBox 42") {
        $tr = $sh.TextFrame.TextRange
        # Stage through an unrelated placeholder string first so the
        # engine's paragraph-text diff doesn't try to reuse/split the run
        # via a common-prefix/suffix match against the old wording; this
        # keeps the edited paragraph a single run (same rPr) like the other
        # untouched paragraphs instead of fragmenting it into two runs.
        $tr.Text = "Bivalve farm`r100 x 480-m longlines with`r80 5-m hang lines spaced 4 m apart`rZZZ_PLACEHOLDER_ZZZ`r"
        $tr.Text = "Bivalve farm`r100 x 480-m longlines with`r80 5-m hang lines spaced 4 m apart`r400 cm of bivalves per foot`r"
        break
    }
}
